$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text (e.g. "43.462.52", "114.26") rather than
# numbers. Force the whole column range to Text format before writing the new
# values so Excel does not silently reinterpret numeric-looking strings as
# numbers, then restore the default ("Normal") style so no stray per-cell style
# index is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.296.55"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.284.94"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "113.72"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "266.77"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "47.92"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "9.36"
$ws.Range("E12").Value = "  +10.15%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "15.54"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "2.616.64"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "0.869"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "2.278.74"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "43.351.81"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").Value = "6.86"
$ws.Range("E20").Value = "  +4.77%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "2.51"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "233.90"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "9.71"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("D27").Value = "11.46"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "40.94"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "173.43"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").Value = "21.49"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "0.0908"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").Value = "5.76"
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "0.0367"
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").Value = "3.96"
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +11.10%  "
$ws.Range("D42").Value = "78.61"
$ws.Range("E42").Value = "  +7.52%  "
$ws.Range("D43").Value = "14.27"
$ws.Range("E43").Value = "  +6.05%  "
$ws.Range("D44").Value = "0.239"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "6.29"
$ws.Range("E45").Value = "  +5.80%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("D48").Value = "8.70"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "104.57"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").Value = "0.0998"
$ws.Range("E51").Value = "  -0.05%  "

$dRange.Style = "Normal"
